$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4716
$ws.Range("H74").Value = 3983.9
$ws.Range("I74").Value = 3990.6
$ws.Range("J74").Value = 3977.2
$ws.Range("K74").Value = 3990.6
$ws.Range("L74").Value = 3977.2
$ws.Range("M74").Value = -3054.6
$ws.Range("N74").Value = -5849.2
$ws.Range("H77").Value = 3983.9
$ws.Range("I77").Value = 3990.6
$ws.Range("J77").Value = 3977.2
$ws.Range("K77").Value = 19953
$ws.Range("L77").Value = 19886
$ws.Range("M77").Value = -15273
$ws.Range("N77").Value = -29246
$ws.Range("H120").Value = 64761
$ws.Range("J120").Value = 64761
$ws.Range("L120").Value = 64761
$ws.Range("N120").Value = -74437
$ws.Range("H138").Value = 3851.2026
$ws.Range("I138").Value = 3413.2856
$ws.Range("J138").Value = 4009.7585
$ws.Range("K138").Value = 10239.8568
$ws.Range("L138").Value = 12029.2755
$ws.Range("M138").Value = -5099.856800000001
$ws.Range("N138").Value = -22309.2755

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4733.0625
$ws.Range("I63").Value = 4111
$ws.Range("J63").Value = 4940.4165
$ws.Range("K63").Value = 4111
$ws.Range("L63").Value = 4940.4165
$ws.Range("M63").Value = -3425
$ws.Range("N63").Value = -6312.4165
$ws.Range("H66").Value = 4733.0625
$ws.Range("I66").Value = 4111
$ws.Range("J66").Value = 4940.4165
$ws.Range("K66").Value = 20555
$ws.Range("L66").Value = 24702.0825
$ws.Range("M66").Value = -17123
$ws.Range("N66").Value = -31566.0825
$ws.Range("H74").Value = 13159697
$ws.Range("I74").Value = 915.36365
$ws.Range("K74").Value = 915.36365
$ws.Range("M74").Value = -41.36365000000001
$ws.Range("H77").Value = 13159697
$ws.Range("I77").Value = 915.36365
$ws.Range("K77").Value = 4576.81825
$ws.Range("M77").Value = -208.8182500000003
$ws.Range("H132").Value = 1285072.4
$ws.Range("I132").Value = 2306.475
$ws.Range("J132").Value = 3850604
$ws.Range("K132").Value = 6919.424999999999
$ws.Range("L132").Value = 11551812
$ws.Range("M132").Value = -4389.424999999999
$ws.Range("N132").Value = -11556872

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2013.875
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2222
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2222
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -4468
$ws.Range("H89").Value = 2013.875
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2222
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 11110
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -22342

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5066.75
$ws.Range("I31").Value = 1388.6
$ws.Range("J31").Value = 7325.263
$ws.Range("K31").Value = 1388.6
$ws.Range("L31").Value = 7325.263
$ws.Range("M31").Value = -1093.6
$ws.Range("N31").Value = -7915.263
$ws.Range("H34").Value = 5066.75
$ws.Range("I34").Value = 1388.6
$ws.Range("J34").Value = 7325.263
$ws.Range("K34").Value = 1388.6
$ws.Range("L34").Value = 7325.263
$ws.Range("M34").Value = -1186.6
$ws.Range("N34").Value = -7729.263
$ws.Range("H45").Value = 1300
$ws.Range("I45").Value = 1300
$ws.Range("K45").Value = 1300
$ws.Range("M45").Value = -707
$ws.Range("H62").Value = 4114.967
$ws.Range("I62").Value = 4042.6428
$ws.Range("J62").Value = 4178.25
$ws.Range("K62").Value = 4042.6428
$ws.Range("L62").Value = 4178.25
$ws.Range("M62").Value = -3418.6428
$ws.Range("N62").Value = -5426.25
$ws.Range("H65").Value = 4114.967
$ws.Range("I65").Value = 4042.6428
$ws.Range("J65").Value = 4178.25
$ws.Range("K65").Value = 20213.214
$ws.Range("L65").Value = 20891.25
$ws.Range("M65").Value = -17093.214
$ws.Range("N65").Value = -27131.25
$ws.Range("H68").Value = 22537.924
$ws.Range("J68").Value = 22537.924
$ws.Range("L68").Value = 22537.924
$ws.Range("N68").Value = -24035.924
$ws.Range("H71").Value = 22537.924
$ws.Range("J71").Value = 22537.924
$ws.Range("L71").Value = 67613.772
$ws.Range("N71").Value = -75101.772

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 689908.9399999999
$ws.Range("J129").Value = 979037.5
$ws.Range("L129").Value = 2937112.5
$ws.Range("N129").Value = -2947112.5
$ws.Range("H138").Value = 2568.524
$ws.Range("I138").Value = 1224.3889
$ws.Range("J138").Value = 10633.333
$ws.Range("K138").Value = 3673.1667
$ws.Range("L138").Value = 31899.999
$ws.Range("M138").Value = 1466.8333
$ws.Range("N138").Value = -42179.999
$ws.Range("H139").Value = 315984.38
$ws.Range("I139").Value = 436533.03
$ws.Range("K139").Value = 1309599.09
$ws.Range("M139").Value = -1304459.09

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 157.91667
$ws.Range("I2").Value = 149.6
$ws.Range("J2").Value = 199.5
$ws.Range("K2").Value = 149.6
$ws.Range("L2").Value = 199.5
$ws.Range("M2").Value = -36.59999999999999
$ws.Range("N2").Value = -425.5
$ws.Range("H70").Value = 5563.684
$ws.Range("I70").Value = 5454.615
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 5454.615
$ws.Range("L70").Value = 5800
$ws.Range("M70").Value = -5184.615
$ws.Range("N70").Value = -6340
$ws.Range("H73").Value = 5563.684
$ws.Range("I73").Value = 5454.615
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 5454.615
$ws.Range("L73").Value = 5800
$ws.Range("M73").Value = -4518.615
$ws.Range("N73").Value = -7672
$ws.Range("H80").Value = 17608846
$ws.Range("I80").Value = 35135110
$ws.Range("J80").Value = 82580
$ws.Range("K80").Value = 35135110
$ws.Range("L80").Value = 82580
$ws.Range("M80").Value = -35134112
$ws.Range("N80").Value = -84576
$ws.Range("H83").Value = 17608846
$ws.Range("I83").Value = 35135110
$ws.Range("J83").Value = 82580
$ws.Range("K83").Value = 175675550
$ws.Range("L83").Value = 412900
$ws.Range("M83").Value = -175670558
$ws.Range("N83").Value = -422884
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 556.5
$ws.Range("I32").Value = 556.5
$ws.Range("K32").Value = 556.5
$ws.Range("M32").Value = -239.5
$ws.Range("H46").Value = 696.6667
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 745
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 745
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1121
$ws.Range("H68").Value = 2962.6667
$ws.Range("J68").Value = 2962.6667
$ws.Range("L68").Value = 2962.6667
$ws.Range("N68").Value = -4460.6667
$ws.Range("H71").Value = 2962.6667
$ws.Range("J71").Value = 2962.6667
$ws.Range("L71").Value = 14813.3335
$ws.Range("N71").Value = -22301.3335
$ws.Range("H100").Value = 73042.62
$ws.Range("I100").Value = 131848.58
$ws.Range("J100").Value = 4435.6665
$ws.Range("K100").Value = 131848.58
$ws.Range("L100").Value = 4435.6665
$ws.Range("M100").Value = -131307.58
$ws.Range("N100").Value = -5517.6665

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4421463
$ws.Range("I132").Value = 2240.8
$ws.Range("K132").Value = 6722.400000000001
$ws.Range("M132").Value = -4192.400000000001
